# Added variability in productivity and quantity of the work: append a new
# column to both the "work_method" sheet (productivity_variation) and the
# "design" sheet (design_variation).

$wb = $excel.ActiveWorkbook

# --- work_method sheet: add "productivity_variation" column (D) ---
$wsWorkMethod = $wb.Worksheets.Item("work_method")
$wsWorkMethod.Range("D1").Value = "productivity_variation"
$wsWorkMethod.Range("D2").Value = 0.5
$wsWorkMethod.Range("D3").Value = 0.5
$wsWorkMethod.Range("D4").Value = 0.5
$wsWorkMethod.Range("D5").Value = 0.5

# --- design sheet: add "design_variation" column (D) ---
$wsDesign = $wb.Worksheets.Item("design")
$wsDesign.Range("D1").Value = "design_variation"
$wsDesign.Range("D2").Value = 3
$wsDesign.Range("D3").Value = 3
$wsDesign.Range("D4").Value = 3
$wsDesign.Range("D5").Value = 3
$wsDesign.Range("D6").Value = 3
$wsDesign.Range("D7").Value = 3

# The new column on "design" nudges column A/D to re-fit their width, same
# as the first column of "work_method"/"quantity" already had (bestFit).
$wsDesign.Columns("A").ColumnWidth = 11.17
$wsDesign.Columns("D").ColumnWidth = 13.45

# Leave the cursor on the newly-added cell in work_method (matches the
# authored selection), then restore "design" as the active sheet/tab.
$wsWorkMethod.Range("D2").Select() | Out-Null
$wsDesign.Activate() | Out-Null
